# Update the "Förändrad" (Changed) date column (C) for rows 2-10,
# incrementing the stored serial date value from 46061 to 46062
# (i.e. advancing the date by one day), matching the automatic
# file-update commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 10; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 46061) {
        $cell.Value2 = 46062
    }
}
